$wb = $excel.ActiveWorkbook

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03d2bada063d4cc94ab60a6499b296d46387686a/e2e/"
$mdFile4df = "4df01f36-eb7a-4837-8ea3-295ed1881c43.md"
$mdFileFff = "ffff3d16f885-6cc2-4bb9-9e3b-459eb75109ac.md"

# ------------------------------------------------------------------
# Overview sheet: handback status text + widened status columns
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ------------------------------------------------------------------
# zh-cn sheet: record the generated handback / target file info
# ------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("J2").Value = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 19:13:14"
$wsZhCn.Range("K3").Value = "2016-08-31 19:13:14"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# rebuild the hyperlinks collection in final row-major order so the
# relationship ids line up the way Excel renumbers them on save
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrlBase + $mdFileFff, [Type]::Missing, [Type]::Missing, $mdFileFff)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)

# ------------------------------------------------------------------
# de-de sheet: record the generated handback / target file info
# ------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("J2").Value = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.de-de.xlf"
$wsDeDe.Range("J3").Value = "4df01f36-eb7a-4837-8ea3-295ed1881c43.096e73f99289f61a7af2c0f49ccedd829c099418.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 19:13:21"
$wsDeDe.Range("K3").Value = "2016-08-31 19:13:21"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrlBase + $mdFileFff, [Type]::Missing, [Type]::Missing, $mdFileFff)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrlBase + $mdFile4df, [Type]::Missing, [Type]::Missing, $mdFile4df)
